$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TopCities")

$ws.Range("A3").Value = "Delhi"
$ws.Range("A4").Value = "Mumbai"
$ws.Range("A5").Value = "Hyderabad"
$ws.Range("A6").Value = "Pune"
